$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values for the updated rows.
# Prices are stored as text (matching the sheet's existing text-formatted Price column),
# so each cell's NumberFormat is forced to Text before the write, then restored to the
# sheet's default (Normal) style so values such as "248.60" or "1.000" are not
# auto-coerced to numbers and lose their trailing/insignificant characters.
$priceUpdates = [ordered]@{
    2 = '30.311.60'
    3 = '1.932.21'
    4 = '1.000'
    5 = '248.60'
    6 = '0.7237'
    8 = '0.3333'
    9 = '28.05'
    10 = '0.06956'
    11 = '0.8084'
    12 = '0.08059'
    13 = '1.932.69'
    14 = '5.435'
    15 = '94.64'
    16 = '14.58'
    17 = '30.304.20'
    18 = '253.68'
    19 = '0.000008202'
    20 = '5.816'
    21 = '2.187.24'
    22 = '0.9999'
    23 = '1.000'
    24 = '6.895'
    25 = '9.759'
    26 = '159.33'
    27 = '2.432'
    29 = '0.1338'
    30 = '1.556'
    31 = '1.337'
    32 = '4.426'
    33 = '4.217'
    34 = '0.05143'
    35 = '1.232'
    36 = '0.7463'
    38 = '0.01996'
    39 = '2.830'
    40 = '6.647'
    41 = '79.00'
    42 = '0.4481'
    43 = '2.003'
    44 = '0.9999'
    45 = '0.8374'
    46 = '102.10'
    47 = '9.784'
    48 = '7.318'
    49 = '36.55'
    50 = '1.494'
    51 = '0.05960'
}

# New "Volume(1h)" (column E) values for the updated rows (plain text, e.g. "  -3.17%  ").
$volumeUpdates = [ordered]@{
    2 = '  -3.17%  '
    3 = '  -3.65%  '
    5 = '  -4.03%  '
    6 = '  -8.48%  '
    7 = '  +0.03%  '
    8 = '  -8.07%  '
    9 = '  -1.47%  '
    10 = '  -1.90%  '
    11 = '  -5.77%  '
    12 = '  -1.60%  '
    13 = '  -3.63%  '
    14 = '  -3.39%  '
    15 = '  -7.00%  '
    16 = '  -3.18%  '
    17 = '  -3.19%  '
    18 = '  -8.16%  '
    19 = '  +2.47%  '
    20 = '  -2.11%  '
    21 = '  -3.41%  '
    22 = '  +0.09%  '
    23 = '  +0.09%  '
    24 = '  -4.27%  '
    25 = '  -3.69%  '
    26 = '  -3.45%  '
    27 = '  +0.89%  '
    28 = '  -4.55%  '
    29 = '  -11.82%  '
    31 = '  -1.58%  '
    32 = '  -4.26%  '
    33 = '  -4.97%  '
    34 = '  -1.87%  '
    35 = '  +0.68%  '
    36 = '  -4.73%  '
    37 = '  -2.30%  '
    38 = '  -0.85%  '
    39 = '  -3.71%  '
    40 = '  -1.19%  '
    41 = '  -2.13%  '
    42 = '  -6.14%  '
    43 = '  -7.44%  '
    44 = '  +0.02%  '
    45 = '  -2.48%  '
    46 = '  -5.27%  '
    47 = '  -1.75%  '
    48 = '  -7.21%  '
    49 = '  -1.17%  '
    50 = '  -0.51%  '
    51 = '  -0.53%  '
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"       # force text so e.g. "248.60"/"1.000" keep their literal form
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"         # restore the cell's original (unstyled) appearance
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
